# remove sy-datum/sy-uzeit from demo programs
#
# Sheet1 originally showed a "Click here to visit abap2xlsx homepage"
# hyperlink plus the current date/time (B2:C4). Sheet2 originally showed a
# "Current Date:" label, the current date, and a "This is link to the
# third sheet" hyperlink (A1:A6). Both are cleared out, leaving blank
# sheets (just the lone placeholder cell that was already on row 1).

$wb = $excel.ActiveWorkbook

# --- Sheet1: drop the homepage hyperlink + date/time block (B2:C4) ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Hyperlinks.Delete()
$ws1.Rows("2:4").Delete()

# --- Sheet2: drop the "third sheet" hyperlink + date label/value block (A1:A6) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Hyperlinks.Delete()
$ws2.Rows("2:6").Delete()
$ws2.Range("A1").ClearContents()
# Touch the style so the now-blank A1 is still written out as a real
# (empty) cell instead of being dropped from sheetData entirely - this
# mirrors Sheet1's row 1, which already carried an empty placeholder cell.
$ws2.Range("A1").Style = "Normal"
